# Update workbook version string after release build.
$wb = $excel.ActiveWorkbook

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $val = $cell.Value2
        if ($val -ne $null -and $val -is [string] -and $val.Contains($oldVersion)) {
            $cell.Value2 = $val.Replace($oldVersion, $newVersion)
        }
    }
}
